$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Work from the LAST paragraph backward to the FIRST so that indices
# of not-yet-processed (lower numbered) paragraphs never shift under
# us, and so that no Delete() call ever spans from an early paragraph
# across a later paragraph boundary into the final paragraph (which
# triggers a bug in this runtime where the delete silently no-ops).
# ------------------------------------------------------------------

# --- Paragraph 13 (last paragraph, keeps the _GoBack bookmark) -----
# Becomes: "Aria also suggests music, to which a user may be
# interested in. ... fluid when operated on a diverse range of
# devices."
$p13 = $d.Paragraphs.Item(13)
$r13 = $d.Range($p13.Range.Start, $p13.Range.End - 1)
$r13.Text = "Aria also suggests music, to which a user may be interested in. At the least, it provides some easy listening. On the other hand it could act as talking point when meeting with their matched musician. Here, the app is trying to extent a user’s musical pallet and influence their playing style with new ideas. Finally it aims to provide a user experience which looks aesthetically pleasing, but fluid when operated on a diverse range of devices."

# --- Paragraph 12 is fully superseded -> delete it whole -----------
$p12 = $d.Paragraphs.Item(12)
$d.Range($p12.Range.Start, $p12.Range.End).Delete()

# --- Paragraph 11 becomes the "On opening the application..." text -
$p11 = $d.Paragraphs.Item(11)
$r11 = $d.Range($p11.Range.Start, $p11.Range.End - 1)
$r11.Text = "On opening the application a user can describe the kind of musician they are looking to work with. Subsequently, they can either leave it for Aria to find someone local, or manually browse through a persona list. On discovering a matching musician, their profile can be saved to a favourites list, to query later when ready. Each profile has a least one link to which a user can follow and contact the matched musician. "

# --- Paragraph 10 ("Sarah" standalone) is removed -------------------
$p10 = $d.Paragraphs.Item(10)
$d.Range($p10.Range.Start, $p10.Range.End).Delete()

# --- Paragraph 9 (empty paragraph) is removed -----------------------
$p9 = $d.Paragraphs.Item(9)
$d.Range($p9.Range.Start, $p9.Range.End).Delete()

# --- Paragraph 8 becomes "When visually separated some musicians..."
$p8 = $d.Paragraphs.Item(8)
$p8.Range.ListFormat.RemoveNumbers()
$p8.Range.Style = "Normal"
$r8 = $d.Range($p8.Range.Start, $p8.Range.End - 1)
$r8.Text = "When visually separated some musicians still have the ability to connect and play in harmony (Schober, 2014). Despite this, we attempt to remove any barrier in the way of two parties joining to exchange musical ideas. Aria facilities the search and connection of musicians who want to collaborate, but can’t find anyone who has similar ambition, ability or musical acuity. "

# --- Paragraph 7 becomes "Brief Overview" ---------------------------
$p7 = $d.Paragraphs.Item(7)
$p7.Range.ListFormat.RemoveNumbers()
$p7.Range.Style = "Normal"
$r7 = $d.Range($p7.Range.Start, $p7.Range.End - 1)
$r7.Text = "Brief Overview"

# --- Paragraph 6 becomes the italic "Little time exceeds..." text --
$p6 = $d.Paragraphs.Item(6)
$p6.Range.ListFormat.RemoveNumbers()
$p6.Range.Style = "Normal"
$r6 = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$r6.Text = "Little time exceeds before she receives a notification from Aria. It turn outs there is bass player only ten minutes from her location, he loves Trip-Hop too! She saves the contact information of the player. Later in the day, Sarah decides to message the player through one of the contact links on his Aria profile. She has now completed her wish to find a bassist who likes similar music. Aria has provided a bridge, now it’s up to both musicians to collaborate."
$p6.Range.Font.Italic = 1

# --- Paragraph 5 becomes the italic "She begins searching..." text -
$p5 = $d.Paragraphs.Item(5)
$p5.Range.ListFormat.RemoveNumbers()
$p5.Range.Style = "Normal"
$r5 = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$r5.Text = "She begins searching through her contacts on social media to see if anyone fits the bill. Sarah quickly realises the search will provide no avail and decides to open Aria on her smartphone. She enters her requirements into the app before leaving home to attend a late lecture. "
$p5.Range.Font.Italic = 1

# --- Paragraph 4 becomes the italic "Sarah, a 20 year old..." text -
$p4 = $d.Paragraphs.Item(4)
$p4.Range.ListFormat.RemoveNumbers()
$p4.Range.Style = "Normal"
$r4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$r4.Text = "Sarah, a 20 year old student who writes music, has an issue; she needs a skilful bass player to accompany her on stage at a gig venue in September. She’s heard from friends that a well-known promoter may be attending. Sarah knows opportunities like this rarely occur. Her psyche yearns to find someone who loves Trip-Hop as much as she does. From experience, she understands the best gigs are played with those of like mind. However Sarah does not know anyone in her circle who could help. "
$p4.Range.Font.Italic = 1

# --- Paragraph 3 becomes plain "Sarah" ------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.ListFormat.RemoveNumbers()
$p3.Range.Style = "Normal"
$r3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$r3.Text = "Sarah"

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
